{"js": "// InforGem: Auto-guardado de Inspecci\u00f3n para el equipo 70-GC-013\n// Applies the set of text edits captured by the commit diff:\n//   1) \"21 de febrero de 2026\" -> \"22 de febrero de 2026\" (6 occurrences)\n//   2) temperature reading 66.5 -> 50.0\n//   3) append \" bien\" to the \"\u00f3ptimas condiciones...\" sentence\n//   4) \"Ignacio\" -> \"Ignacio perro\" (2 occurrences)\n//   5) \"Pendiente\" -> \"Emi\u00e1n nama\" (1 occurrence)\n\nconst body = context.document.body;\n\n// 1) Date: 21 de febrero de 2026 -> 22 de febrero de 2026 (appears 6x,\n//    both as a standalone cell/paragraph value and inline inside longer\n//    sentences such as \"Firmado por Lorena Rojas el 21 de febrero de 2026\").\nconst dateMatches = body.search(\"21 de febrero de 2026\", { matchCase: true });\ndateMatches.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < dateMatches.items.length; i++) {\n  dateMatches.items[i].insertText(\"22 de febrero de 2026\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Temperature reading changes from 66.5 to 50.0 inside the sentence.\nconst tempMatches = body.search(\n  \"Verificaci\u00f3n de par\u00e1metros de operaci\u00f3n (Presi\u00f3n de carga:   / descarga:  y temperatura de salida del elemento (66.5).\",\n  { matchCase: true }\n);\ntempMatches.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < tempMatches.items.length; i++) {\n  tempMatches.items[i].insertText(\n    \"Verificaci\u00f3n de par\u00e1metros de operaci\u00f3n (Presi\u00f3n de carga:   / descarga:  y temperatura de salida del elemento (50.0).\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// 3) Append \" bien\" to the equipment status sentence.\nconst condMatches = body.search(\n  \"El equipo se encuentra funcionando en \u00f3ptimas condiciones...\",\n  { matchCase: true }\n);\ncondMatches.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < condMatches.items.length; i++) {\n  condMatches.items[i].insertText(\n    \"El equipo se encuentra funcionando en \u00f3ptimas condiciones... bien\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// 4) \"Ignacio\" -> \"Ignacio perro\" (the technician name cell + the\n//    \"Firmado por Ignacio\" run \u2014 2 occurrences total).\nconst ignacioMatches = body.search(\"Ignacio\", { matchCase: true, matchWholeWord: true });\nignacioMatches.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < ignacioMatches.items.length; i++) {\n  ignacioMatches.items[i].insertText(\"Ignacio perro\", \"Replace\");\n}\nawait context.sync();\n\n// 5) \"Pendiente\" -> \"Emi\u00e1n nama\" (technician name placeholder in the table).\nconst pendienteMatches = body.search(\"Pendiente\", { matchCase: true, matchWholeWord: true });\npendienteMatches.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < pendienteMatches.items.length; i++) {\n  pendienteMatches.items[i].insertText(\"Emi\u00e1n nama\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# InforGem: Auto-guardado de Inspecci\u00f3n para el equipo 70-GC-013\n# Applies the set of text edits captured by the commit diff:\n#   1) \"21 de febrero de 2026\" -> \"22 de febrero de 2026\" (6 occurrences)\n#   2) temperature reading 66.5 -> 50.0\n#   3) append \" bien\" to the \"\u00f3ptimas condiciones...\" sentence\n#   4) \"Ignacio\" -> \"Ignacio perro\" (2 occurrences)\n#   5) \"Pendiente\" -> \"Emi\u00e1n nama\" (1 occurrence)\n\n$d = $word.ActiveDocument\n\n# 1) Date: 21 de febrero de 2026 -> 22 de febrero de 2026 (6 occurrences,\n#    both standalone cell/paragraph values and inline inside longer\n#    sentences such as \"Firmado por Lorena Rojas el 21 de febrero de 2026\").\n$findDate = $d.Content.Find\n$findDate.Execute(\"21 de febrero de 2026\", $false, $false, $false, $false, $false, $true, 1, $false, \"22 de febrero de 2026\", 2)\n\n# 2) Temperature reading changes from 66.5 to 50.0 inside the sentence.\n$findTemp = $d.Content.Find\n$findTemp.Execute(\"temperatura de salida del elemento (66.5).\", $true, $false, $false, $false, $false, $true, 1, $false, \"temperatura de salida del elemento (50.0).\", 2)\n\n# 3) Append \" bien\" to the equipment status sentence.\n$findCond = $d.Content.Find\n$findCond.Execute(\"El equipo se encuentra funcionando en \u00f3ptimas condiciones...\", $true, $false, $false, $false, $false, $true, 1, $false, \"El equipo se encuentra funcionando en \u00f3ptimas condiciones... bien\", 2)\n\n# 4) \"Ignacio\" -> \"Ignacio perro\" (the technician name cell + the\n#    \"Firmado por Ignacio\" run - 2 occurrences total).\n$findIgnacio = $d.Content.Find\n$findIgnacio.Execute(\"Ignacio\", $true, $true, $false, $false, $false, $true, 1, $false, \"Ignacio perro\", 2)\n\n# 5) \"Pendiente\" -> \"Emi\u00e1n nama\" (technician name placeholder in the table).\n$findPendiente = $d.Content.Find\n$findPendiente.Execute(\"Pendiente\", $true, $true, $false, $false, $false, $true, 1, $false, \"Emi\u00e1n nama\", 2)\n"}
